# Gaussian Quadrature export: rename the sheet tab to the shorter "UniformF"
# title and append the extra quadrature-scheme row (row 16) that mirrors the
# existing "HexGrid-60degTilt5degRes" row (row 15), but indexed 14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet tab rename: "UniformF-HW50.xpc" -> "UniformF"
$ws.Name = "UniformF"

# New row 16 reuses row 15's formatting (bold/border/centered "s=1" style
# on column A) via copy/paste-format, then gets its own values.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A16").Value = 14

# Column B repeats the same label text as row 15 (shared string reused).
$ws.Range("B16").Value = $ws.Range("B15").Value()

# Columns C:M are all 1, same as every other data row.
$ws.Range("C16:M16").Value = 1

Write-Host "Row 16 added and sheet renamed"
